$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 21:50"

# Update country labels that shifted position (column A)
$ws.Range("A73").Value = "Camerun"
$ws.Range("A74").Value = "Libano"
$ws.Range("A75").Value = "Tunez"
$ws.Range("A76").Value = "Letonia"
$ws.Range("A77").Value = "Bulgaria"
$ws.Range("A78").Value = "Kazajistan"
$ws.Range("A79").Value = "Eslovaquia"
$ws.Range("A80").Value = "Azerbaiyan"
$ws.Range("A81").Value = "Principado de Andorra"
$ws.Range("A82").Value = "Republica de Macedonia"
$ws.Range("A83").Value = "Kuwait"
$ws.Range("A84").Value = "Costa Rica"
$ws.Range("A85").Value = "Republica de Chipre"
$ws.Range("A86").Value = "Uruguay"
$ws.Range("A87").Value = "Bielorrusia"
$ws.Range("A88").Value = "Taiwan"
$ws.Range("A89").Value = "Reunion"
$ws.Range("A90").Value = "Jordania"
$ws.Range("A143").Value = "Puerto Rico"
$ws.Range("A144").Value = "Mali"
$ws.Range("A145").Value = "Zambia"
$ws.Range("A151").Value = "Guyana"
$ws.Range("A152").Value = "San Martin (Parte Holandesa)"
$ws.Range("A153").Value = "Eritrea"
$ws.Range("A154").Value = "San Martin (Parte Francesa)"
$ws.Range("A155").Value = "Congo"
$ws.Range("A156").Value = "Gabon"
$ws.Range("A157").Value = "Birmania"
$ws.Range("A158").Value = "Tanzania"
$ws.Range("A160").Value = "Nueva Caledonia"
$ws.Range("A161").Value = "Haiti"
$ws.Range("A164").Value = "Siria"
$ws.Range("A165").Value = "Benin"

# Update numeric data cells (columns B-H)
$ws.Range("B4").Value = 269996
$ws.Range("C4").Value = 25119
$ws.Range("D4").Value = 12015
$ws.Range("E4").Value = 251057
$ws.Range("F4").Value = 5787
$ws.Range("G4").Value = 854
$ws.Range("H4").Value = 6924
$ws.Range("F15").Value = 1324
$ws.Range("D16").Value = 2186
$ws.Range("E16").Value = 9981
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = 208
$ws.Range("B35").Value = 2684
$ws.Range("C35").Value = 263
$ws.Range("E35").Value = 2518
$ws.Range("E43").Value = 1473
$ws.Range("G43").Value = 9
$ws.Range("H43").Value = 62
$ws.Range("E58").Value = 1023
$ws.Range("G58").Value = 5
$ws.Range("H58").Value = 27
$ws.Range("B73").Value = 509
$ws.Range("C73").Value = 203
$ws.Range("D73").Value = 17
$ws.Range("E73").Value = 484
$ws.Range("F73").Value = 0
$ws.Range("H73").Value = 8
$ws.Range("B74").Value = 508
$ws.Range("C74").Value = 14
$ws.Range("D74").Value = 50
$ws.Range("E74").Value = 441
$ws.Range("F74").Value = 26
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 17
$ws.Range("B75").Value = 495
$ws.Range("C75").Value = 40
$ws.Range("D75").Value = 5
$ws.Range("E75").Value = 472
$ws.Range("F75").Value = 30
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 18
$ws.Range("B76").Value = 493
$ws.Range("C76").Value = 35
$ws.Range("D76").Value = 1
$ws.Range("E76").Value = 491
$ws.Range("F76").Value = 3
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 1
$ws.Range("B77").Value = 485
$ws.Range("C77").Value = 28
$ws.Range("D77").Value = 30
$ws.Range("E77").Value = 441
$ws.Range("F77").Value = 18
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 14
$ws.Range("B78").Value = 464
$ws.Range("C78").Value = 29
$ws.Range("D78").Value = 29
$ws.Range("E78").Value = 429
$ws.Range("F78").Value = 6
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 6
$ws.Range("B79").Value = 450
$ws.Range("C79").Value = 24
$ws.Range("D79").Value = 10
$ws.Range("E79").Value = 439
$ws.Range("F79").Value = 3
$ws.Range("H79").Value = 1
$ws.Range("B80").Value = 443
$ws.Range("C80").Value = 43
$ws.Range("D80").Value = 32
$ws.Range("E80").Value = 406
$ws.Range("F80").Value = 7
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 5
$ws.Range("B81").Value = 439
$ws.Range("C81").Value = 11
$ws.Range("D81").Value = 16
$ws.Range("E81").Value = 407
$ws.Range("F81").Value = 12
$ws.Range("H81").Value = 16
$ws.Range("B82").Value = 430
$ws.Range("C82").Value = 46
$ws.Range("D82").Value = 20
$ws.Range("E82").Value = 398
$ws.Range("F82").Value = 8
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 12
$ws.Range("B83").Value = 417
$ws.Range("C83").Value = 75
$ws.Range("D83").Value = 82
$ws.Range("E83").Value = 335
$ws.Range("F83").Value = 16
$ws.Range("H83").Value = 0
$ws.Range("B84").Value = 416
$ws.Range("C84").Value = 20
$ws.Range("D84").Value = 11
$ws.Range("E84").Value = 403
$ws.Range("F84").Value = 13
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 2
$ws.Range("B85").Value = 396
$ws.Range("C85").Value = 40
$ws.Range("D85").Value = 28
$ws.Range("E85").Value = 357
$ws.Range("F85").Value = 11
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 11
$ws.Range("B86").Value = 369
$ws.Range("C86").Value = 19
$ws.Range("D86").Value = 68
$ws.Range("E86").Value = 297
$ws.Range("F86").Value = 13
$ws.Range("B87").Value = 351
$ws.Range("C87").Value = 47
$ws.Range("D87").Value = 53
$ws.Range("E87").Value = 294
$ws.Range("F87").Value = 11
$ws.Range("H87").Value = 4
$ws.Range("B88").Value = 348
$ws.Range("C88").Value = 9
$ws.Range("D88").Value = 50
$ws.Range("E88").Value = 293
$ws.Range("F88").Value = 0
$ws.Range("H88").Value = 5
$ws.Range("B89").Value = 321
$ws.Range("C89").Value = 13
$ws.Range("D89").Value = 40
$ws.Range("E89").Value = 281
$ws.Range("F89").Value = 3
$ws.Range("H89").Value = 0
$ws.Range("B90").Value = 310
$ws.Range("C90").Value = 11
$ws.Range("D90").Value = 58
$ws.Range("E90").Value = 247
$ws.Range("F90").Value = 5
$ws.Range("H90").Value = 5
$ws.Range("D143").Value = 1
$ws.Range("H143").Value = 2
$ws.Range("C144").Value = 3
$ws.Range("D144").Value = 0
$ws.Range("H144").Value = 3
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 2
$ws.Range("H145").Value = 1
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 19
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 4
$ws.Range("B152").Value = 23
$ws.Range("C152").Value = 5
$ws.Range("D152").Value = 6
$ws.Range("E152").Value = 15
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 2
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 22
$ws.Range("H153").Value = 0
$ws.Range("E154").Value = 19
$ws.Range("H154").Value = 1
$ws.Range("B155").Value = 22
$ws.Range("D155").Value = 2
$ws.Range("E155").Value = 18
$ws.Range("H155").Value = 2
$ws.Range("B156").Value = 21
$ws.Range("D156").Value = 1
$ws.Range("D157").Value = 0
$ws.Range("E157").Value = 19
$ws.Range("B158").Value = 20
$ws.Range("D158").Value = 3
$ws.Range("E158").Value = 16
$ws.Range("H158").Value = 1
$ws.Range("C160").Value = 0
$ws.Range("C161").Value = 2
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 0
$ws.Range("H164").Value = 2
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 2
$ws.Range("H165").Value = 0
